$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new attendance row (row 3), mirroring the existing row 2 pattern:
# Date (text), EmployeeID (number), ClockIn (text), ClockOut (blank text), Log (blank text)
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Range("A3").Value = "10/05/2025"
$ws.Range("B3").Value = 2602069620
$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Range("C3").Value = "14:05:37"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
